$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-12 Saturday", "2025-04-13 Sunday"),
    @("50÷7=7, 1", "90÷2=45, 0"),
    @("10÷6=1, 4", "69÷8=8, 5"),
    @("11÷2=5, 1", "16÷7=2, 2"),
    @("36÷5=7, 1", "49÷3=16, 1"),
    @("87÷7=12, 3", "84÷5=16, 4"),
    @("63÷6=10, 3", "77÷5=15, 2"),
    @("77÷4=19, 1", "26÷7=3, 5"),
    @("74÷4=18, 2", "24÷7=3, 3"),
    @("73÷2=36, 1", "58÷3=19, 1"),
    @("30÷3=10, 0", "88÷9=9, 7"),
    @("45÷2=22, 1", "30÷5=6, 0"),
    @("91÷8=11, 3", "99÷5=19, 4"),
    @("93÷5=18, 3", "14÷9=1, 5"),
    @("58÷2=29, 0", "39÷6=6, 3"),
    @("72÷3=24, 0", "67÷8=8, 3"),
    @("65÷7=9, 2", "72÷5=14, 2"),
    @("39÷3=13, 0", "71÷9=7, 8"),
    @("51÷4=12, 3", "80÷7=11, 3"),
    @("88÷5=17, 3", "33÷3=11, 0"),
    @("66÷8=8, 2", "32÷2=16, 0"),
    @("13÷8=1, 5", "47÷3=15, 2"),
    @("96÷2=48, 0", "63÷4=15, 3"),
    @("95÷2=47, 1", "84÷9=9, 3"),
    @("62÷7=8, 6", "62÷8=7, 6"),
    @("18÷4=4, 2", "51÷3=17, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND:" $old
    }
}

$d.Save()
